$d = $word.ActiveDocument

# 1) Remove the "_GoBack" bookmark (bookmarkStart/bookmarkEnd pair) left over
#    from the previous editing session. "_GoBack" is a hidden bookmark so it
#    does not show up in the default Bookmarks collection/Count, but it can
#    still be reached (and deleted) by name.
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    if ($goBack -ne $null) {
        $goBack.Delete()
    }
} catch {
}

# 2) The ">>>  your stuff after this line >>>" paragraph currently has
#    grammar-check proofing artefacts splitting it into three runs.
#    Re-asserting the same text collapses it back into a single run.
$d.Content.Find.Execute(">>>  your stuff after this line >>>", $true, $false, $false, $false, $false, $true, 1, $false, ">>>  your stuff after this line >>>", 2)

# 3) "Ben changing things up!" -> "Zhen Ma changing things up!"
$d.Content.Find.Execute("Ben changing things up!", $true, $false, $false, $false, $false, $true, 1, $false, "Zhen Ma changing things up!", 2)

# 4) Replace the first of the two trailing blank paragraphs (the one right
#    after "...changing things up!") with the new "Zoey create a rfp!"
#    sentence, keeping the spell-check proofing marks Word places around the
#    flagged word "rfp". The second blank paragraph is left untouched.
$target = $null
for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "`r" -and $i -gt 1) {
        $prev = $d.Paragraphs($i - 1)
        if ($prev.Range.Text -like "*changing things up!*") {
            $target = $para
            break
        }
    }
}

if ($target -eq $null) {
    # Fallback: first blank paragraph in the document.
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs($i)
        if ($para.Range.Text -eq "`r") {
            $target = $para
            break
        }
    }
}

if ($target -ne $null) {
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:r><w:t xml:space="preserve">Zoey create a </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:r><w:t>rfp</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t>!</w:t></w:r>' +
           '</w:p>'
    $target.Range.InsertXML($xml)
}
